$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (rows 2..6) for columns D, K, L, M, N, O, P, R, S
# These values come from a cyclic re-ordering of the original rows.

$data = @(
    @{ Row = 2; D = 44301; K = "Hachiya"; L = "Segunda"; M = 250; N = 20000; O = 21000; P = 20500; R = "Región de O'Higgins"; S = 1139 },
    @{ Row = 3; D = 44355; K = "Mankaki"; L = "Segunda"; M = 270; N = 20000; O = 21000; P = 20500; R = "Región Metropolitana"; S = 1139 },
    @{ Row = 4; D = 44342; K = "Mankaki"; L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 },
    @{ Row = 5; D = 44305; K = "Mankaki"; L = "Segunda"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 },
    @{ Row = 6; D = 44313; K = "Mankaki"; L = "Primera"; M = 270; N = 21000; O = 22000; P = 21500; R = "Región de O'Higgins"; S = 1194 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("D$r").Value = $item.D
    $ws.Range("K$r").Value = $item.K
    $ws.Range("L$r").Value = $item.L
    $ws.Range("M$r").Value = $item.M
    $ws.Range("N$r").Value = $item.N
    $ws.Range("O$r").Value = $item.O
    $ws.Range("P$r").Value = $item.P
    $ws.Range("R$r").Value = $item.R
    $ws.Range("S$r").Value = $item.S
}
